$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Define new text content (order matters: it controls shared-string index assignment)
$s_B9 = @'
A few of many agent deskop features are: 
1. See the full picture with contextual details about customer history and preferences
2.  Manage all conversations in one streamlined interface. Field inbound messages or place outbound calls using a personalized directory
3. View and track individual agent performance metrics like AHT, FCR, and more
'@
$s_B11 = @'
Free package ($0/Month) offers:
1. Automated Routing 
2.Interactive Voice Response
3. Chatbots
4. Analytics and Reporting
'@
$s_B10 = @'
Mega Cloud is available in three flavors at a variety of affordable price points. 
1. Free - $0/mo
2. Premium - $99/mo
3. Platinum - $170/mo
'@
$s_B12 = @'
Premium package ($99/Month) offers:
1. Automated Routing
2.Interactive Voice Response 
3. Chatbots
4. Analytics and Reporting
5. Third Party Integration
'@
$s_B13 = @'
Platinum package ($170/Month) offers:
1. Automated Routing
2.Interactive Voice Response
3. Chatbots 
4. Analytics and Reporting
5. Third Party Integration
'@
$s_B14 = @'
Free package ($0/Month) offers:
1. Automated Routing
2.Interactive Voice Response
3. Chatbots
4. Analytics and Reporting 
Premium package ($99/Month) offers: 
1. Automated Routing 
2.Interactive Voice Response
3. Chatbots
4. Analytics and Reporting
5. Third Party Integration
Platinum package ($170/Month) offers: 
1. Automated Routing 
2. Interactive Voice Response 
3. Chatbots
4. Analytics and Reporting
5. Third Party Integration
'@
$s_B2_run2 = @'
  provides personalized, omnichannel customer support from a single interface through our cloud-based contact center solution. Some of the attractive features of Mega Cloud are:
1. Omnichannel Management
2. Future-Proofed Contact Center
3.Chatbots with Text-to-Speech Capabilties
4.Workforce Management Integration
5.Built-In AI Routing Two-Way Messaging
'@

# Row 9: "What are the agent desktop features?" answer
$ws.Range("B9").Value = $s_B9

# Row 11: "What features does Free package offer?" answer
$ws.Range("B11").Value = $s_B11
$ws.Rows(11).RowHeight = 72

# Row 10: "What Mega Cloud packages are available" answer
$ws.Range("B10").Value = $s_B10
$ws.Rows(10).RowHeight = 57.6

# Row 12: "What features does Premium package offer?" answer
$ws.Range("B12").Value = $s_B12
$ws.Rows(12).RowHeight = 86.4

# Row 13: "What features does Platinum package offer?" answer
$ws.Range("B13").Value = $s_B13
$ws.Rows(13).RowHeight = 86.4

# Row 14: "What are the features of each available package?" answer (combined)
$ws.Range("B14").Value = $s_B14
$ws.Rows(14).RowHeight = 273.6

# Row 2: "What features does Mega Cloud offer?" answer - rich text w/ bold "Mega Cloud" lead-in
$ws.Range("B2").Value = "Mega Cloud" + $s_B2_run2
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Characters(11, $s_B2_run2.Length).Font.Bold = $false
$ws.Range("B2").Characters(1, 10).Font.Bold = $true
$ws.Range("B2").Font.Bold = $false

# Update view: select B3 and scroll sheet back to the top-left
$ws.Activate()
$ws.Range("B3").Select()
